$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet ---
$ws.Name = "Sidwell"

# --- Delete rows that are no longer needed ---
# Remove one duplicate blank filler row and 6 mid-table data rows so the sheet
# shrinks from 34 rows to 27 rows; old rows 32/33/34 (blank/sum/spacer) become
# new rows 25/26/27, preserving their formatting.
$ws.Rows("31:31").Delete()
$ws.Rows("8:13").Delete()

# --- Title row ---
$ws.Range("A1").Value = 'ELEN-7046 Group Project – Individual Time sheet - Sidwell Mokhemisa'

# --- Data rows ---
$ws.Range("A3").Value = 'Project Group Meeting 1'
$ws.Range("B3").Value = 42484
$ws.Range("C3").Value = 0.6041666666666666
$ws.Range("D3").Value = 0.6875
$ws.Range("E3").Value = 2

$ws.Range("A4").Value = 'Start High Level Design Activities'
$ws.Range("B4").Value = 42488
$ws.Range("C4").Value = 0.375
$ws.Range("D4").Value = 0.7083333333333334
$ws.Range("E4").Value = 8

$ws.Range("A5").Value = 'Project Group Meeting 2'
$ws.Range("B5").Value = 42491
$ws.Range("C5").Value = 0.6041666666666666
$ws.Range("D5").Value = 0.6875
$ws.Range("E5").Value = 2

$ws.Range("A6").Value = 'Start Use Case Modeling'
$ws.Range("B6").Value = 42492
$ws.Range("C6").Value = 0.7291666666666666
$ws.Range("D6").Value = 0.8958333333333334
$ws.Range("E6").Value = 10

$ws.Range("A7").Value = 'Complete Use Case Modeling'
$ws.Range("B7").Value = 42497
$ws.Range("C7").Value = 0.4166666666666667
$ws.Range("D7").Value = 0.625
$ws.Range("E7").Value = 8

$ws.Range("A8").Value = 'Project Group Meeting 3'
$ws.Range("B8").Value = 42498
$ws.Range("C8").Value = 0.6041666666666666
$ws.Range("D8").Value = 0.6875
$ws.Range("E8").Value = 2

$ws.Range("A9").Value = 'Identify Delivery Method (SDLC)'
$ws.Range("B9").Value = 42504
$ws.Range("C9").Value = 0.4166666666666667
$ws.Range("D9").Value = 0.5416666666666666
$ws.Range("E9").Value = 4

$ws.Range("A10").Value = 'Project Group Meeting 4'
$ws.Range("B10").Value = 42505
$ws.Range("C10").Value = 0.6041666666666666
$ws.Range("D10").Value = 0.6875
$ws.Range("E10").Value = 2

$ws.Range("A11").Value = 'Deliver a Tailored Method based on RUP'
$ws.Range("B11").Value = 42506
$ws.Range("C11").Value = 0.6666666666666666
$ws.Range("D11").Value = 0.9166666666666666
$ws.Range("E11").Value = 3

$ws.Range("A12").Value = 'Project Group Meeting 5'
$ws.Range("B12").Value = 42512
$ws.Range("C12").Value = 0.6041666666666666
$ws.Range("D12").Value = 0.6875
$ws.Range("E12").Value = 2

$ws.Range("A13").Value = 'Component Model'
$ws.Range("B13").Value = 42513
$ws.Range("C13").Value = 0.75
$ws.Range("D13").Value = 0.875
$ws.Range("E13").Value = 9

$ws.Range("A14").Value = 'Project Group Meeting 6'
$ws.Range("B14").Value = 42526
$ws.Range("C14").Value = 0.5
$ws.Range("D14").Value = 0.75
$ws.Range("E14").Value = 2

$ws.Range("A15").Value = 'Infrastructure Deign'
$ws.Range("B15").Value = 42527
$ws.Range("C15").Value = 0.75
$ws.Range("D15").Value = 0.9166666666666666
$ws.Range("E15").Value = 6

$ws.Range("A16").Value = 'Project Group Meeting 7'
$ws.Range("B16").Value = 42532
$ws.Range("C16").Value = 0.4166666666666667
$ws.Range("D16").Value = 0.7916666666666666
$ws.Range("E16").Value = 2

$ws.Range("A17").Value = 'Project Group Meeting 8'
$ws.Range("B17").Value = 42539
$ws.Range("C17").Value = 0.375
$ws.Range("D17").Value = 0.875
$ws.Range("E17").Value = 8

$ws.Range("A18").Value = 'Project Group Meeting 9'
$ws.Range("B18").Value = 42540
$ws.Range("C18").Value = 0.375
$ws.Range("D18").Value = 0.5416666666666666
$ws.Range("E18").Value = 2

$ws.Range("A19").Value = 'Project Presentation Deck'
$ws.Range("B19").Value = 42540
$ws.Range("C19").Value = 0.5833333333333334
$ws.Range("D19").Value = 0.75
$ws.Range("E19").Value = 3

$ws.Range("A20").Value = 'Project Group Meeting 10'
$ws.Range("B20").Value = 42547
$ws.Range("C20").Value = 0.6041666666666666
$ws.Range("D20").Value = 0.7291666666666666
$ws.Range("E20").Value = 2

$ws.Range("A21").Value = 'Reports - Individual and Group'
$ws.Range("B21").Value = 42551
$ws.Range("C21").Value = 0.375
$ws.Range("D21").Value = 0.7708333333333334
$ws.Range("E21").Value = 16

$ws.Range("A22").Value = 'Group Meeting 11(Virtual)'
$ws.Range("B22").Value = 42551
$ws.Range("C22").Value = 0.7708333333333334
$ws.Range("D22").Value = 0.8333333333333334
$ws.Range("E22").Value = 1

$ws.Range("A23").Value = 'Project Group Meeting 12'
$ws.Range("B23").Value = 42552
$ws.Range("C23").Value = 0.4166666666666667
$ws.Range("D23").Value = 0.7083333333333334
$ws.Range("E23").Value = 8

$ws.Range("A24").Value = 'Individual Report'
$ws.Range("B24").Value = 42553
$ws.Range("C24").Value = 0.4166666666666667
$ws.Range("D24").Value = 0.9166666666666666
$ws.Range("E24").Value = 11

# --- Break F6 out of the shared-formula group (kept as its own literal formula) ---
$ws.Range("F6").Formula = "=D6-C6"

# --- Clear the blank row leftovers (row 25) ---
$ws.Range("A25:F25").ClearContents()

$wb.Save()
